$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$ws.Range("T5").Value = [double]"5.503011245186134E-5"
$ws.Range("U5").Value = [double]"8.3540650792029026E-4"
$ws.Range("T6").Value = [double]"5.9967451738382029E-5"
$ws.Range("U6").Value = [double]"3.7823060819885398E-4"
$ws.Range("T7").Value = [double]"5.2591582908256337E-5"
$ws.Range("U7").Value = [double]"6.1545257336610053E-4"
$ws.Range("T8").Value = [double]"6.3097607068156441E-5"
$ws.Range("U8").Value = [double]"3.1850768399081672E-4"
$ws.Range("T9").Value = [double]"6.651413183117218E-5"
$ws.Range("U9").Value = [double]"3.9108209196853691E-4"
$ws.Range("T10").Value = [double]"5.6136414158981842E-5"
$ws.Range("U10").Value = [double]"4.6807925904620537E-4"
$ws.Range("T11").Value = [double]"5.9176816500790682E-5"
$ws.Range("U11").Value = [double]"4.2690692371642752E-4"
$ws.Range("T12").Value = [double]"5.3486824799421879E-5"
$ws.Range("U12").Value = [double]"6.2754703219099824E-4"
$ws.Range("T13").Value = [double]"6.4052016992818749E-5"
$ws.Range("U13").Value = [double]"3.6169760037931009E-4"
$ws.Range("T14").Value = [double]"7.244650571687623E-5"
$ws.Range("U14").Value = [double]"3.3460575962597901E-4"
$ws.Range("T15").Value = [double]"5.8788226974038367E-5"
$ws.Range("U15").Value = [double]"5.3558070035386682E-4"
$ws.Range("T16").Value = [double]"6.4415858830368177E-5"
$ws.Range("U16").Value = [double]"5.0167872385648601E-4"
$ws.Range("T17").Value = [double]"7.1043242439360123E-5"
$ws.Range("U17").Value = [double]"3.0788156582625838E-4"
$ws.Range("T18").Value = [double]"6.2696668505210143E-5"
$ws.Range("U18").Value = [double]"3.2990491081884331E-4"
$ws.Range("T19").Value = [double]"6.587384557238129E-5"
$ws.Range("U19").Value = [double]"4.4467908621744271E-4"
$ws.Range("T20").Value = [double]"6.8778554865638961E-5"
$ws.Range("U20").Value = [double]"3.2966716429677341E-4"
$ws.Range("T21").Value = [double]"6.1696713670308831E-5"
$ws.Range("U21").Value = [double]"6.1832939730611892E-4"
$ws.Range("T22").Value = [double]"6.2271385162116001E-5"
$ws.Range("U22").Value = [double]"4.9750038377983051E-4"
$ws.Range("T23").Value = [double]"6.6534327743237738E-5"
$ws.Range("U23").Value = [double]"3.9604336967912737E-4"
$ws.Range("T24").Value = [double]"6.7953094678194638E-5"
$ws.Range("U24").Value = [double]"2.9339056320201859E-4"
$ws.Range("T25").Value = [double]"6.3636899634327115E-5"
$ws.Range("U25").Value = [double]"3.6736937705477199E-4"
$ws.Range("T26").Value = [double]"5.5042517614221591E-5"
$ws.Range("U26").Value = [double]"4.8190988349849212E-4"
$ws.Range("T27").Value = [double]"6.1508612493463589E-5"
$ws.Range("U27").Value = [double]"3.730815843846586E-4"
$ws.Range("T28").Value = [double]"6.5575908474247891E-5"
$ws.Range("U28").Value = [double]"4.3861615287973032E-4"
$ws.Range("T29").Value = [double]"5.3472650915101847E-5"
$ws.Range("U29").Value = [double]"6.6416436890334306E-4"
$ws.Range("T30").Value = [double]"5.589026740771058E-5"
$ws.Range("U30").Value = [double]"5.6239866346945719E-4"
$ws.Range("T31").Value = [double]"6.0667109582902133E-5"
$ws.Range("U31").Value = [double]"3.4072890349385432E-4"
$ws.Range("T32").Value = [double]"6.4819447104103053E-5"
$ws.Range("U32").Value = [double]"5.1475647798412971E-4"
$ws.Range("T33").Value = [double]"5.7980135711774843E-5"
$ws.Range("U33").Value = [double]"6.6343539150521451E-4"
$ws.Range("T34").Value = [double]"6.2409254788756301E-5"
$ws.Range("U34").Value = [double]"4.5863123904253602E-4"
$ws.Range("T35").Value = [double]"5.7053680281262648E-5"
$ws.Range("U35").Value = [double]"5.9131996851791613E-4"
$ws.Range("T36").Value = [double]"5.9313111471295708E-5"
$ws.Range("U36").Value = [double]"3.0412926905250672E-4"
$ws.Range("T37").Value = [double]"5.8873679508800583E-5"
$ws.Range("U37").Value = [double]"4.9599222719775053E-4"
$ws.Range("T38").Value = [double]"6.7396970948967267E-5"
$ws.Range("U38").Value = [double]"3.0102327399908512E-4"
$ws.Range("T39").Value = [double]"6.6479926343178161E-5"
$ws.Range("U39").Value = [double]"4.2070986772552578E-4"
$ws.Range("T40").Value = [double]"5.7016136681753607E-5"
$ws.Range("U40").Value = [double]"4.3018509130313381E-4"
$ws.Range("T41").Value = [double]"6.2959941453997841E-5"
$ws.Range("U41").Value = [double]"7.8432089532466838E-4"
$ws.Range("T42").Value = [double]"5.3312671823951763E-5"
$ws.Range("U42").Value = [double]"5.9893724041657501E-4"
$ws.Range("T43").Value = [double]"6.0662792959156002E-5"
$ws.Range("U43").Value = [double]"5.3786702840837529E-4"
$ws.Range("T44").Value = [double]"6.0029534792445907E-5"
$ws.Range("U44").Value = [double]"6.6301868090244014E-4"
$ws.Range("T45").Value = [double]"6.7094581281543304E-5"
$ws.Range("U45").Value = [double]"3.9966658594246812E-4"
$ws.Range("T46").Value = [double]"6.5098897339210227E-5"
$ws.Range("U46").Value = [double]"4.0487964608311103E-4"
$ws.Range("T47").Value = [double]"5.9981150841951513E-5"
$ws.Range("U47").Value = [double]"5.1727586868207873E-4"
$ws.Range("T48").Value = [double]"6.1855361938343333E-5"
$ws.Range("U48").Value = [double]"3.6377727648352062E-4"
$ws.Range("T49").Value = [double]"6.7824709597586656E-5"
$ws.Range("U49").Value = [double]"5.1329421531259357E-4"
$ws.Range("T50").Value = [double]"5.4776555583699901E-5"
$ws.Range("U50").Value = [double]"5.9226057513606469E-4"
$ws.Range("T51").Value = [double]"6.5138314692705338E-5"
$ws.Range("U51").Value = [double]"2.8794104521623012E-4"
$ws.Range("T52").Value = [double]"6.2254630953650577E-5"
$ws.Range("U52").Value = [double]"3.8433863290953388E-4"
$ws.Range("T53").Value = [double]"6.4450193365712664E-5"
$ws.Range("U53").Value = [double]"3.3611032757183942E-4"
$ws.Range("T54").Value = [double]"6.1309046924711992E-5"
$ws.Range("U54").Value = [double]"3.363964903392386E-4"

# Scroll the view so row 22 is at the top-left, then land the selection on W52
# (mirrors the sheetView's topLeftCell="A22" / selection activeCell="W52" in the target).
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("W52").Select()
